$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: step 1 now fails with an error ---
$ws.Range("L2").Value = "FAIL"
$ws.Range("M2").Value = "Cannot read properties of null (reading 'isClosed')"
$ws.Range("N2").Value = "Cannot read properties of null (reading 'isClosed')"
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""

# --- Row 3: step not executed anymore (TO BE EXECUTED -> NO), result columns cleared ---
$ws.Range("A3").Value = "NO"
$ws.Range("L3:P3").ClearContents()

# --- Row 4: step not executed anymore (TO BE EXECUTED -> NO), result columns cleared ---
$ws.Range("A4").Value = "NO"
$ws.Range("L4:P4").ClearContents()

# --- Row 5 & 6: not executed anymore ---
$ws.Range("A5").Value = "NO"
$ws.Range("A6").Value = "NO"

# --- Row 9: step description corrected ---
$ws.Range("G9").Value = "Click New"

# --- Rows 34-40: Flipkart scenario not executed anymore ---
$ws.Range("A34").Value = "NO"
$ws.Range("A35").Value = "NO"
$ws.Range("A36").Value = "NO"
$ws.Range("A37").Value = "NO"
$ws.Range("A38").Value = "NO"
$ws.Range("A39").Value = "NO"
$ws.Range("A40").Value = "NO"
